# "Carregando dados de out/2025" -- re-saving the monthly export.
# The workbook was re-opened/re-saved from a PT-BR desktop Excel: the
# selection moved down two rows, the page margins were reset to the
# locale's metric "Normal" margins (1.3 / 2 / 0.8 cm), and the theme's
# accent1/accent5 swatches came back swapped relative to the prior save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Selection now sits on the two blank rows right after the data
# (previously it spanned three rows) -- also clears the old scrolled
# top-left / zoom state from the view.
$null = $ws.Range("A155:XFD156").Select()

# Page margins: 1.3cm / 2cm / 0.8cm (metric "Normal" margins), expressed
# in points (1 inch = 72 pt, 1 inch = 2.54 cm) for the COM PageSetup API.
$ws.PageSetup.LeftMargin = 36.850393728
$ws.PageSetup.RightMargin = 36.850393728
$ws.PageSetup.TopMargin = 56.692913399999995
$ws.PageSetup.BottomMargin = 56.692913399999995
$ws.PageSetup.HeaderMargin = 22.67716464
$ws.PageSetup.FooterMargin = 22.67716464

# Theme accent1/accent5 swatches come back swapped (5B9BD5 <-> 4472C4).
$tcs = $wb.Theme.ThemeColorScheme
$tcs.Colors(5).RGB = 13998939   # accent1 -> 5B9BD5
$tcs.Colors(9).RGB = 12874308   # accent5 -> 4472C4
